$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 2
    3  = -6
    4  = 3
    5  = -1
    6  = -2
    7  = -4
    8  = -1
    9  = 1
    10 = -1
    11 = -5
    12 = -1
    13 = -4
    14 = -1
    15 = -2
    16 = 3
    17 = -1
    19 = -3
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
